$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $ws.Columns("A").ColumnWidth
Write-Host $ws.Columns("B").ColumnWidth
